# Commit: "binary search for problem 719"
# Add a new entry to the "新题" (new problems) sheet: solved LeetCode 719
# via binary search on 2019-03-31 (serial date 43555, same day as the
# three rows already logged above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# New row is row 30 (right after the existing last row, 29).
$newRow = 30
$prevRow = $newRow - 1

# Copy the date cell's formatting from the row above so the new date cell
# keeps the same date number-format style instead of creating a new one.
$ws.Cells.Item($prevRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null

# A30: date (same day as the previous couple of rows)
$ws.Cells.Item($newRow, 1).Value = 43555
# B30: problem number
$ws.Cells.Item($newRow, 2).Value = 719
# F30: status
$ws.Cells.Item($newRow, 6).Value = "done"

# Reflect the author's final view state on that sheet.
$ws.Activate()
$ws.Range("K27").Select() | Out-Null
